# Regenerate orders with updated distance/size codes.
#
# The experiment's distance conditions were renumbered
#   D51 -> D55, D64 -> D69, D80 -> D86
# and the "large" size condition was renumbered
#   S30 -> S31
#
# These tokens show up all over the sheet: in the Condition column
# (e.g. "Face15_D51_S20"), in the Filename_Left / Filename_Right columns
# (e.g. "Face15_D51_S20_l.png", "Fixation_D80_r.png"), and in the
# standalone Distance / Size lookup columns ("D51", "S30", ...). Every
# occurrence is a plain substring of the cell's text, so walk every used
# cell and rewrite any string value in place.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$rows = $used.Rows.Count
$cols = $used.Columns.Count

for ($r = 1; $r -le $rows; $r++) {
  for ($c = 1; $c -le $cols; $c++) {
    $cell = $ws.Cells.Item($r, $c)
    $v = $cell.Value2
    if ($v -is [string]) {
      $nv = $v.Replace("D51", "D55").Replace("D64", "D69").Replace("D80", "D86").Replace("S30", "S31")
      if ($nv -ne $v) {
        $cell.Value = $nv
      }
    }
  }
}
